$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-spelled NPC name "HellsGateKeeper" -> "Hell'sGatekeeper"
# (this is the `name` column value for NPC id 13, row 14)
$ws.Range("B14").Value = "Hell'sGatekeeper"

# Fixes to weekly fight x/y position values for NPC id 13 (row 14)
$ws.Range("I14").Value = 1552
$ws.Range("J14").Value = 1392

# Fixes to weekly fight x/y position values for NPC id 33 (row 34)
$ws.Range("I34").Value = 592
$ws.Range("J34").Value = 336
